# Updated cryptos list on Wed Oct 25 16:46:11 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($Row, $Price, $Volume) {
    $priceCell = $ws.Cells.Item($Row, 4)
    $priceCell.NumberFormat = "@"
    $priceCell.Value = $Price
    $ws.Cells.Item($Row, 5).Value = $Volume
}

function Set-FullRow($Row, $Coin, $Link, $Price, $Volume) {
    $ws.Cells.Item($Row, 2).Value = $Coin
    $ws.Cells.Item($Row, 3).Value = $Link
    Set-Row $Row $Price $Volume
}

Set-Row 2  "34.685.09"  "  +1.81%  "
Set-Row 3  "1.792.78"   "  +0.26%  "
Set-Row 5  "225.09"     "  -0.04%  "
Set-Row 6  "0.553"      "  -0.91%  "
Set-Row 7  "1.00"       "  +0.50%  "
Set-Row 8  "32.71"      "  +5.77%  "
Set-Row 9  "0.284"      "  +1.60%  "
Set-Row 10 "0.0668"     "  +1.02%  "
Set-Row 11 "0.0937"     "  +1.45%  "
Set-Row 12 "2.052.66"   "  +0.40%  "

Set-FullRow 13 "WrappedEther" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth" "1.813.76" "  +1.51%  "
Set-FullRow 14 "Chainlink"    "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"    "11.07"    "  +10.80%  "

Set-Row 15 "0.635"      "  +1.14%  "
Set-Row 16 "34.649.46"  "  +1.92%  "
Set-Row 17 "4.28"       "  +1.89%  "
Set-Row 18 "69.12"      "  +0.83%  "
Set-Row 19 "254.17"     "  +0.87%  "
Set-Row 20 "0.0₃0762"   "  +2.94%  "
Set-Row 21 "0.999"      "  +0.40%  "
Set-Row 22 "10.37"      "  +0.65%  "
Set-Row 23 "4.21"       "  -0.13%  "
Set-Row 24 "2.13"       "  -1.10%  "
Set-Row 25 "158.59"     "  +0.67%  "
Set-Row 26 "16.38"      "  -0.58%  "
Set-Row 27 "7.08"       "  +1.99%  "
Set-Row 29 "1.00"       "  +0.55%  "
Set-Row 30 "0.0518"     "  +1.29%  "
Set-Row 31 "3.76"       "  -1.23%  "
Set-Row 32 "1.19"       "  -0.25%  "
Set-Row 33 "3.57"       "  +1.13%  "
Set-Row 34 "1.86"       "  +6.25%  "
Set-Row 35 "1.448.30"   "  -3.13%  "
Set-Row 36 "1.06"       "  -0.23%  "
Set-Row 37 "0.0189"     "  +1.60%  "
Set-Row 38 "0.625"      "  -0.54%  "
Set-Row 39 "82.97"      "  -0.14%  "
Set-Row 40 "2.82"       "  +4.02%  "
Set-Row 41 "2.35"       "  -0.11%  "
Set-Row 42 "0.899"      "  +1.00%  "
Set-Row 43 "2.07"       "  -0.81%  "
Set-Row 44 "0.0507"     "  -0.51%  "
Set-Row 45 "5.94"       "  +3.34%  "

Set-FullRow 46 "WEMIXToken"    "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"   "1.05"      "  -1.38%  "
Set-FullRow 47 "RocketPoolETH" "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth" "1.952.88"  "  +0.62%  "

Set-Row 48 "104.54"     "  +6.50%  "
Set-Row 49 "1.00"       "  +0.36%  "
Set-Row 50 "11.89"      "  -0.64%  "

Set-FullRow 51 "BabyDogeCoin" "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge" "0.0₆0122" "  +4.62%  "
